$wb = $excel.ActiveWorkbook

# Generate Report for handback: the zh-cn and de-de sheets each gain a
# "Latest Target File" (E) and "Latest Handback File" (F) hyperlink for
# the two source files, the Status (B) flips from "Ready for handoff" to
# "Handed back: in sync with en-US", and the Latest Handback DateTime (G)
# records the actual handback time.

$rows = @(
  @{ Sheet = "zh-cn"; Row = 2; MdName = "3bb7d2ea-d3e9-4813-8fbe-f10d3b6343df.md";
     MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/7c3b4417aeb34ec936a5bc3f751328da611b49b0/e2e/3bb7d2ea-d3e9-4813-8fbe-f10d3b6343df.md";
     XlfName = "3bb7d2ea-d3e9-4813-8fbe-f10d3b6343df.29bc73d6262f16dd1c6d2e4d704324e64ac0a4cf.zh-cn.xlf";
     XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/48a4f65ccc73d06c7aad6dd37e11639dec234298/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/3bb7d2ea-d3e9-4813-8fbe-f10d3b6343df.29bc73d6262f16dd1c6d2e4d704324e64ac0a4cf.zh-cn.xlf";
     HandbackTime = "2016-01-18 04:58:52" },
  @{ Sheet = "zh-cn"; Row = 3; MdName = "7ca93064-f50f-4eee-8df3-5a021638a516.md";
     MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/7c3b4417aeb34ec936a5bc3f751328da611b49b0/e2e/7ca93064-f50f-4eee-8df3-5a021638a516.md";
     XlfName = "7ca93064-f50f-4eee-8df3-5a021638a516.20eb648f70b7ded1c03d50847477a172d8524a17.zh-cn.xlf";
     XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/48a4f65ccc73d06c7aad6dd37e11639dec234298/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/7ca93064-f50f-4eee-8df3-5a021638a516.20eb648f70b7ded1c03d50847477a172d8524a17.zh-cn.xlf";
     HandbackTime = "2016-01-18 04:58:52" },
  @{ Sheet = "de-de"; Row = 2; MdName = "3bb7d2ea-d3e9-4813-8fbe-f10d3b6343df.md";
     MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/7c3b4417aeb34ec936a5bc3f751328da611b49b0/e2e/3bb7d2ea-d3e9-4813-8fbe-f10d3b6343df.md";
     XlfName = "3bb7d2ea-d3e9-4813-8fbe-f10d3b6343df.29bc73d6262f16dd1c6d2e4d704324e64ac0a4cf.de-de.xlf";
     XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/15abd386bc8cba4d4b05ffc901a2f52e4b35fff6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/3bb7d2ea-d3e9-4813-8fbe-f10d3b6343df.29bc73d6262f16dd1c6d2e4d704324e64ac0a4cf.de-de.xlf";
     HandbackTime = "2016-01-18 04:59:14" },
  @{ Sheet = "de-de"; Row = 3; MdName = "7ca93064-f50f-4eee-8df3-5a021638a516.md";
     MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/7c3b4417aeb34ec936a5bc3f751328da611b49b0/e2e/7ca93064-f50f-4eee-8df3-5a021638a516.md";
     XlfName = "7ca93064-f50f-4eee-8df3-5a021638a516.20eb648f70b7ded1c03d50847477a172d8524a17.de-de.xlf";
     XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/15abd386bc8cba4d4b05ffc901a2f52e4b35fff6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/7ca93064-f50f-4eee-8df3-5a021638a516.20eb648f70b7ded1c03d50847477a172d8524a17.de-de.xlf";
     HandbackTime = "2016-01-18 04:59:14" }
)

foreach ($r in $rows) {
    $ws = $wb.Worksheets.Item($r.Sheet)

    # Status -> handed back, now in sync with en-US
    $ws.Cells.Item($r.Row, 2).Value = "Handed back: in sync with en-US"

    # Latest Target File (E) - link to the source markdown, same as column A
    $eCell = $ws.Cells.Item($r.Row, 5)
    $ws.Hyperlinks.Add($eCell, $r.MdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $r.MdName)

    # Latest Handback File (F) - link to the handed-back xlf, same as column C
    $fCell = $ws.Cells.Item($r.Row, 6)
    $ws.Hyperlinks.Add($fCell, $r.XlfUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $r.XlfName)

    # Latest Handback DateTime (G)
    $ws.Cells.Item($r.Row, 7).Value = $r.HandbackTime
}
